# Apply new calculated statistics values to both worksheets (Chan1, Chan2).
# Commit message: "new method for cal - excel engine use the old xlwt, it may cause problems"
# This reflects a recomputation of the statistics values in rows 2 and 3
# (row 1 is the header row and is left untouched) on both sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (Chan1) ---
$ws1 = $wb.Worksheets.Item("Chan1")

$ws1.Range("A2").Value = 994
$ws1.Range("B2").Value = 0.3329979879275654
$ws1.Range("C2").Value = 0.05533199195171026
$ws1.Range("D2").Value = 0.6116700201207244
$ws1.Range("E2").Value = 0.1167002012072435
$ws1.Range("F2").Value = 0.02414486921529175
$ws1.Range("G2").Value = 0.0482897384305835
$ws1.Range("H2").Value = 0.04426559356136821

$ws1.Range("A3").Value = 1319
$ws1.Range("B3").Value = 0.2759666413949962
$ws1.Range("C3").Value = 0.2615617892342684
$ws1.Range("D3").Value = 0.4624715693707354
$ws1.Range("E3").Value = 0.1902956785443518
$ws1.Range("F3").Value = 0.03411675511751327
$ws1.Range("G3").Value = 0.1053828658074299
$ws1.Range("H3").Value = 0.05079605761940864

# --- Sheet 2 (Chan2) ---
$ws2 = $wb.Worksheets.Item("Chan2")

$ws2.Range("A2").Value = 1267
$ws2.Range("B2").Value = 0.2541436464088398
$ws2.Range("C2").Value = 0.1846882399368587
$ws2.Range("D2").Value = 0.5611681136543015
$ws2.Range("E2").Value = 0.1104972375690608
$ws2.Range("F2").Value = 0.009471191791633781
$ws2.Range("G2").Value = 0.06945540647198106
$ws2.Range("H2").Value = 0.03157063930544594

$ws2.Range("A3").Value = 1655
$ws2.Range("B3").Value = 0.2066465256797583
$ws2.Range("C3").Value = 0.04531722054380664
$ws2.Range("D3").Value = 0.748036253776435
$ws2.Range("E3").Value = 0.1087613293051359
$ws2.Range("F3").Value = 0.02114803625377644
$ws2.Range("G3").Value = 0.0338368580060423
$ws2.Range("H3").Value = 0.05377643504531722

$wb.Save()
